# ---------------------------------------------------------------------------
# Applies the "stroke" outcome sheet + data fixes described by the diff:
#  - fixes the "Dagnas" -> "Dangas" author typo (all sheets that use it)
#  - adds a new A2 value "CRAP STUDY" (row 2 trial name) on every outcome sheet
#  - appends a new study row (row 12, "4D-ACS" / Jang 2025 / prasugrel) to every
#    outcome sheet
#  - adds a brand new "stroke" outcome worksheet with the same study rows
#  - updates sheet selections / the active-tab bookmark to match the diff
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # all_cause_mortality
$ws2 = $wb.Worksheets.Item(2)   # major_or_relevant_bleeding
$ws3 = $wb.Worksheets.Item(3)   # mi

# ---------------------------------------------------------------------------
# Phase 1: fix the "Dagnas" -> "Dangas" typo everywhere it is used (B4 on each
# of the three existing sheets). Doing this first - and on every sheet that
# references it - lets the old shared string drop out and "Dangas" becomes
# the first newly-created shared string.
# ---------------------------------------------------------------------------
$ws1.Range("B4").Value = "Dangas"
$ws2.Range("B4").Value = "Dangas"
$ws3.Range("B4").Value = "Dangas"

# ---------------------------------------------------------------------------
# Phase 2: append the new "4D-ACS" / Jang study row (row 12) to sheet1, in an
# order that creates new shared strings as "Jang", "4D-ACS", "prasugrel"
# (matching the target shared-string table order).
# ---------------------------------------------------------------------------
$ws1.Range("B12").Value = "Jang"
$ws1.Range("A12").Value = "4D-ACS"
$ws1.Range("D12").Value = "prasugrel"
$ws1.Range("C12").Value = 2025
$ws1.Range("E12").Value = 5
$ws1.Range("F12").Value = 328
$ws1.Range("G12").Value = 2
$ws1.Range("H12").Value = 328
$ws1.Range("B11:H11").Copy()
$ws1.Range("B12:H12").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Phase 3: add the new A2 "CRAP STUDY" trial name to sheet1 (last new shared
# string created).
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = "CRAP STUDY"

# ---------------------------------------------------------------------------
# Phase 4: replicate the same row-12 + A2 additions onto sheet2 and sheet3.
# All strings already exist at this point, so no further shared strings are
# created - the remaining sheets simply reuse them.
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "CRAP STUDY"
$ws2.Range("A12").Value = "4D-ACS"
$ws2.Range("B12").Value = "Jang"
$ws2.Range("C12").Value = 2025
$ws2.Range("D12").Value = "prasugrel"
$ws2.Range("E12").Value = 2
$ws2.Range("F12").Value = 328
$ws2.Range("G12").Value = 15
$ws2.Range("H12").Value = 328
$ws2.Range("B11:H11").Copy()
$ws2.Range("B12:H12").PasteSpecial(-4122)

$ws3.Range("A2").Value = "CRAP STUDY"
$ws3.Range("A12").Value = "4D-ACS"
$ws3.Range("B12").Value = "Jang"
$ws3.Range("C12").Value = 2025
$ws3.Range("D12").Value = "prasugrel"
$ws3.Range("E12").Value = 1
$ws3.Range("F12").Value = 328
$ws3.Range("G12").Value = 3
$ws3.Range("H12").Value = 328
$ws3.Range("B11:H11").Copy()
$ws3.Range("B12:H12").PasteSpecial(-4122)
# sheet3's new A12 picked up style "2" in the original edit (copied from the
# header-row formatting) - replicate that quirk.
$ws3.Range("D1").Copy()
$ws3.Range("A12").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Phase 5: create the new "stroke" worksheet by duplicating the now fully
# populated "mi" sheet (so every style / shared-string reference lines up
# automatically), then overwrite its event counts with the stroke-specific
# numbers.
# ---------------------------------------------------------------------------
$ws3.Copy([System.Reflection.Missing]::Value, $ws3)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "stroke"

$ws4.Range("E2").Value = 3
$ws4.Range("F2").Value = 826
$ws4.Range("G2").Value = 7
$ws4.Range("H2").Value = 854

$ws4.Range("E3").Value = 23
$ws4.Range("F3").Value = 2283
$ws4.Range("G3").Value = 13
$ws4.Range("H3").Value = 2287

$ws4.Range("E4").Value = 1
$ws4.Range("F4").Value = 1158
$ws4.Range("G4").Value = 2
$ws4.Range("H4").Value = 1184

$ws4.Range("E5").Value = 0
$ws4.Range("F5").Value = 260
$ws4.Range("G5").Value = 3
$ws4.Range("H5").Value = 238

$ws4.Range("E6").Value = 0
$ws4.Range("F6").Value = 245
$ws4.Range("G6").Value = 3
$ws4.Range("H6").Value = 274

$ws4.Range("E7").Value = 4
$ws4.Range("F7").Value = 735
$ws4.Range("G7").Value = 4
$ws4.Range("H7").Value = 738

$ws4.Range("E8").Value = 3
$ws4.Range("F8").Value = 588
$ws4.Range("G8").Value = 3
$ws4.Range("H8").Value = 608

$ws4.Range("E9").Value = 20
$ws4.Range("F9").Value = 1700
$ws4.Range("G9").Value = 24
$ws4.Range("H9").Value = 1700

$ws4.Range("E10").Value = 20
$ws4.Range("F10").Value = 1712
$ws4.Range("G10").Value = 15
$ws4.Range("H10").Value = 1698

$ws4.Range("E11").Value = 3
$ws4.Range("F11").Value = 961
$ws4.Range("G11").Value = 2
$ws4.Range("H11").Value = 981

$ws4.Range("E12").Value = 0
$ws4.Range("F12").Value = 328
$ws4.Range("G12").Value = 4
$ws4.Range("H12").Value = 328

# ---------------------------------------------------------------------------
# Phase 6: sheet-view / selection bookkeeping to match the diff - the active
# tab moves from sheet1 to sheet3 ("mi"), and each sheet's stored selection
# changes.
# ---------------------------------------------------------------------------
$ws1.Range("A2").Select()
$ws2.Range("A2").Select()
$ws4.Range("A2").Select()

$ws3.Range("D23").Select()
$excel.ActiveWindow.DisplayGridlines = $excel.ActiveWindow.DisplayGridlines

$ws3.Activate()

$excel.ActiveWindow.WindowState = $excel.ActiveWindow.WindowState
